# Arreglo Datos de proveedor FT. Agrego Robo FT y no FT.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers for the "Robo" (theft) fast-track block ---
$ws.Range("Y1").Value = "ROBO"
$ws.Range("Z1").Value = "DescripcionRobo"
$ws.Range("AA1").Value = "CodProveedorRobo"

# --- Fix proveedor / poliza / fecha data for rows 2-5 (NroPoliza + FechaSiniestro) ---
# Leading apostrophe keeps these as literal text (matching the existing
# quotePrefix / text-like cell styles) instead of Excel auto-coercing the
# numeric-looking policy number or the dd/mm/yyyy string into a date serial.
$ws.Range("F2").Value = "'04104018548"
$ws.Range("H2").Value = "'08/08/2021"

$ws.Range("F3").Value = "'04104018548"
$ws.Range("H3").Value = "'08/08/2021"

$ws.Range("F4").Value = "'04104018548"
$ws.Range("H4").Value = "'08/08/2021"
$ws.Range("V4").Value = ""

$ws.Range("F5").Value = "'04104018548"
$ws.Range("H5").Value = "'08/08/2021"

# CodProveedor fix on row 8 (was 23766, should be 27433)
$ws.Range("V8").Value = "'27433"

# --- Robo (theft) FastTrack columns Y/Z/AA per row ---
$ws.Range("Y2").Value = "No"
$ws.Range("Y3").Value = "No"

$ws.Range("Y4").Value = "Sí"
$ws.Range("Z4").Value = "Rueda"
$ws.Range("AA4").Value = 27433

$ws.Range("Y5").Value = "Sí"
$ws.Range("Z5").Value = "Batería"

$ws.Range("Y6").Value = "No"
$ws.Range("Y7").Value = "No"
$ws.Range("Y8").Value = "No"
$ws.Range("Y9").Value = "No"

# EsFastTrack (T) was blank for rows 10-13; diff fills it with "No"
$ws.Range("T10").Value = "No"
$ws.Range("Y10").Value = "No"

$ws.Range("T11").Value = "No"
$ws.Range("Y11").Value = "No"

$ws.Range("T12").Value = "No"
$ws.Range("Y12").Value = "No"

$ws.Range("T13").Value = "No"
$ws.Range("Y13").Value = "No"

$ws.Range("Y14").Value = "No"
$ws.Range("Y15").Value = "No"
$ws.Range("Y16").Value = "No"
$ws.Range("Y17").Value = "No"

# --- View state: match the author's scroll position / active selection ---
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H7").Select()
